$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy header style from an existing header cell (A1) onto the new headers
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) # xlPasteFormats

# Fill team record for every data row (2-45)
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 29).Value = 77
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 0
}
